$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 250.47368
$ws.Range("J2").Value = 242.16667
$ws.Range("L2").Value = 242.16667
$ws.Range("N2").Value = -468.16667
$ws.Range("H12").Value = 1098.5
$ws.Range("I12").Value = 130.33333
$ws.Range("J12").Value = 2066.6667
$ws.Range("K12").Value = 130.33333
$ws.Range("L12").Value = 2066.6667
$ws.Range("M12").Value = 39.66667000000001
$ws.Range("N12").Value = -2406.6667
$ws.Range("H38").Value = 6210.9653
$ws.Range("J38").Value = 6638.4443
$ws.Range("L38").Value = 19915.3329
$ws.Range("N38").Value = -20659.3329
$ws.Range("H98").Value = 1494.6364
$ws.Range("I98").Value = 1621.2222
$ws.Range("J98").Value = 925
$ws.Range("K98").Value = 1621.2222
$ws.Range("L98").Value = 925
$ws.Range("M98").Value = -123.2221999999999
$ws.Range("N98").Value = -3921
$ws.Range("H112").Value = 5666.3335
$ws.Range("I112").Value = 1999.5
$ws.Range("K112").Value = 5998.5
$ws.Range("M112").Value = -4890.5
$ws.Range("H115").Value = 1785.0834
$ws.Range("I115").Value = 1202.1
$ws.Range("K115").Value = 3606.3
$ws.Range("M115").Value = -2039.3
$ws.Range("H116").Value = 4249.9165
$ws.Range("I116").Value = 4200
$ws.Range("J116").Value = 4285.5713
$ws.Range("K116").Value = 4200
$ws.Range("L116").Value = 4285.5713
$ws.Range("M116").Value = -758
$ws.Range("N116").Value = -11169.5713
$ws.Range("H122").Value = 1494.6364
$ws.Range("I122").Value = 1621.2222
$ws.Range("J122").Value = 925
$ws.Range("K122").Value = 4863.6666
$ws.Range("L122").Value = 2775
$ws.Range("M122").Value = -2413.6666
$ws.Range("N122").Value = -7675
$ws.Range("H125").Value = 42439.2
$ws.Range("I125").Value = 1997.5
$ws.Range("J125").Value = 69400.336
$ws.Range("K125").Value = 17977.5
$ws.Range("L125").Value = 624603.024
$ws.Range("M125").Value = -15517.5
$ws.Range("N125").Value = -629523.024
$ws.Range("H137").Value = 1919.0857
$ws.Range("I137").Value = 1902.0312
$ws.Range("K137").Value = 5706.0936
$ws.Range("M137").Value = -3156.0936

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5820.645
$ws.Range("I32").Value = 4967.6924
$ws.Range("J32").Value = 10256
$ws.Range("K32").Value = 4967.6924
$ws.Range("L32").Value = 10256
$ws.Range("M32").Value = -4680.6924
$ws.Range("N32").Value = -10830
$ws.Range("H82").Value = 34036
$ws.Range("J82").Value = 37045
$ws.Range("L82").Value = 37045
$ws.Range("N82").Value = -37767
$ws.Range("H85").Value = 34036
$ws.Range("J85").Value = 37045
$ws.Range("L85").Value = 37045
$ws.Range("N85").Value = -39541
$ws.Range("H97").Value = 476.8421
$ws.Range("I97").Value = 462.13333
$ws.Range("J97").Value = 532
$ws.Range("K97").Value = 462.13333
$ws.Range("L97").Value = 532
$ws.Range("M97").Value = 33.86667
$ws.Range("N97").Value = -1524
$ws.Range("H122").Value = 3790.647
$ws.Range("I122").Value = 3703.4167
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 11110.2501
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -8660.250100000001
$ws.Range("N122").Value = -16900
$ws.Range("H132").Value = 1533.6274
$ws.Range("I132").Value = 1508.8511
$ws.Range("K132").Value = 4526.5533
$ws.Range("M132").Value = -1996.5533

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 6886.4614
$ws.Range("I105").Value = 6335
$ws.Range("K105").Value = 6335
$ws.Range("M105").Value = -4588
$ws.Range("H119").Value = 59855
$ws.Range("J119").Value = 59855
$ws.Range("L119").Value = 59855
$ws.Range("N119").Value = -69531

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H33").Value = 2571.4285
$ws.Range("I33").Value = 3200
$ws.Range("K33").Value = 3200
$ws.Range("M33").Value = -2821
$ws.Range("H58").Value = 2430.4443
$ws.Range("I58").Value = 2339.1428
$ws.Range("K58").Value = 2339.1428
$ws.Range("M58").Value = -2136.1428
$ws.Range("H132").Value = 1249.1395
$ws.Range("I132").Value = 1249.1395
$ws.Range("K132").Value = 3747.4185
$ws.Range("M132").Value = -1217.4185
$ws.Range("H134").Value = 2124.5806
$ws.Range("I134").Value = 2050.9167
$ws.Range("K134").Value = 6152.750100000001
$ws.Range("M134").Value = -3617.750100000001
$ws.Range("H136").Value = 2430.4443
$ws.Range("I136").Value = 2339.1428
$ws.Range("K136").Value = 7017.428400000001
$ws.Range("M136").Value = -4467.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 2706.5293
$ws.Range("I18").Value = 1037.4546
$ws.Range("K18").Value = 3112.3638
$ws.Range("M18").Value = -2943.3638
$ws.Range("H109").Value = 13209.4
$ws.Range("I109").Value = 349
$ws.Range("J109").Value = 32500
$ws.Range("K109").Value = 1047
$ws.Range("L109").Value = 97500
$ws.Range("M109").Value = -7
$ws.Range("N109").Value = -99580
$ws.Range("H113").Value = 3300.3
$ws.Range("I113").Value = 1500
$ws.Range("J113").Value = 3500.3333
$ws.Range("K113").Value = 4500
$ws.Range("L113").Value = 10500.9999
$ws.Range("M113").Value = -2330
$ws.Range("N113").Value = -14840.9999
$ws.Range("H125").Value = 15000
$ws.Range("J125").Value = 15000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H138").Value = 3748.75
$ws.Range("J138").Value = 3996
$ws.Range("L138").Value = 11988
$ws.Range("N138").Value = -22268

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 49333
$ws.Range("J64").Value = 49333
$ws.Range("L64").Value = 49333
$ws.Range("N64").Value = -49829
$ws.Range("H67").Value = 49333
$ws.Range("J67").Value = 49333
$ws.Range("L67").Value = 49333
$ws.Range("N67").Value = -51049
$ws.Range("H80").Value = 4553.55
$ws.Range("I80").Value = 3148
$ws.Range("K80").Value = 3148
$ws.Range("M80").Value = -2150
$ws.Range("H83").Value = 4553.55
$ws.Range("I83").Value = 3148
$ws.Range("K83").Value = 15740
$ws.Range("M83").Value = -10748
$ws.Range("H97").Value = 774.7727
$ws.Range("I97").Value = 716.4761999999999
$ws.Range("J97").Value = 1999
$ws.Range("K97").Value = 716.4761999999999
$ws.Range("L97").Value = 1999
$ws.Range("M97").Value = -220.4761999999999
$ws.Range("N97").Value = -2991
$ws.Range("H108").Value = 41770.168
$ws.Range("J108").Value = 40000
$ws.Range("L108").Value = 40000
$ws.Range("N108").Value = -47680
$ws.Range("H122").Value = 1950.5
$ws.Range("I122").Value = 1901
$ws.Range("K122").Value = 5703
$ws.Range("M122").Value = -3253

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1397.4375
$ws.Range("I16").Value = 1238.5
$ws.Range("K16").Value = 1238.5
$ws.Range("M16").Value = -1068.5
$ws.Range("H46").Value = 830.8333
$ws.Range("H55").Value = 1230.0385
$ws.Range("I55").Value = 349.33334
$ws.Range("J55").Value = 3211.625
$ws.Range("K55").Value = 349.33334
$ws.Range("L55").Value = 3211.625
$ws.Range("M55").Value = -176.33334
$ws.Range("N55").Value = -3557.625

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2088.111
$ws.Range("I96").Value = 1368.2307
$ws.Range("K96").Value = 1368.2307
$ws.Range("M96").Value = 4.76929999999993
$ws.Range("H120").Value = 173666.67
$ws.Range("J120").Value = 173666.67
$ws.Range("L120").Value = 173666.67
$ws.Range("N120").Value = -183342.67
